$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")
$ws.Activate()

# Insert a new row at position 83, pushing the blank separator row and the
# "Stunden insgesamt" summary row down to 84/85.
$ws.Rows.Item(83).Insert()

# Carry over the date/number/time formatting used by the row above (row 82)
# for the columns that need it (F..K) before filling in the new entry.
$ws.Range("F82:K82").Copy()
$ws.Range("F83:K83").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new "LaTex Workshop" entry.
$ws.Cells.Item(83, 1).Value = 17
$ws.Cells.Item(83, 2).Value = "Interface Design"
$ws.Cells.Item(83, 3).Value = "[SEMINAR]"
$ws.Cells.Item(83, 4).Value = "LaTex Workshop"
$ws.Cells.Item(83, 5).Value = "How to LaTex"
$ws.Cells.Item(83, 6).Value = (Get-Date -Year 2021 -Month 6 -Day 8).Date
$ws.Cells.Item(83, 7).Value = (Get-Date -Year 2021 -Month 6 -Day 12).Date
$ws.Cells.Item(83, 8).Formula = "=ROUNDUP(((SUM(K83-J83)*24*60/60)/0.25),0)*0.25"
$ws.Cells.Item(83, 10).Value = 0.41666666666666669
$ws.Cells.Item(83, 11).Value = 0.53125

# Match the author's final selection/view state.
$ws.Range("K83").Select()
